# B6-PowerPoint.pptx edit
#
# 1) Three tables (slides 14, 15, 16) get a new "quick style" applied from
#    the Table Styles gallery - this swaps the tableStyleId from the
#    document's custom Table_0 style ({756F1139-...}) to a built-in
#    PowerPoint gallery style ({CAC69970-...}).
# 2) The deck's design/theme colours are switched from the "Integral"
#    (Red Violet) palette to the plain "Office" palette.

$p = $ppt.ActivePresentation

$newTableStyle = "{CAC69970-249F-4958-9B4B-790541CFD7EB}"
$tableSlides = @(14, 15, 16)

foreach ($slideIdx in $tableSlides) {
    $slide = $p.Slides.Item($slideIdx)
    $shape = $slide.Shapes.Item(1)
    if ($shape.HasTable) {
        $shape.Table.ApplyStyle($newTableStyle)
    }
}

# Re-colour the presentation's theme (affects the slide master, and so
# every slide built from it) from the "Integral" / Red Violet scheme to
# the plain "Office" colour scheme.
$officeColors = @(
    0,         # dk1      000000
    16777215,  # lt1      FFFFFF
    6968388,   # dk2      44546A
    15132391,  # lt2      E7E6E6
    13998939,  # accent1  5B9BD5
    3243501,   # accent2  ED7D31
    10855845,  # accent3  A5A5A5
    49407,     # accent4  FFC000
    12874308,  # accent5  4472C4
    4697456,   # accent6  70AD47
    12673797,  # hlink    0563C1
    7491477    # folHlink 954F72
)

$master = $p.SlideMaster
$themeColors = $master.Theme.ThemeColorScheme
for ($i = 1; $i -le $themeColors.Count; $i++) {
    $themeColors.Item($i).RGB = $officeColors[$i - 1]
}
